$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (same shape/columns as existing table), values taken
# from the re-scraped data (same matches, re-ordered).
$rows = @(
    @(" Dubai (DSC)", " October 13 2020", "Super Kings won by 20 runs", "Chennai Super Kings", "Sunrisers Hyderabad", "MS Dhoni (c)†", "21", "13", "2", "1", "161.53"),
    @(" Abu Dhabi", " September 19 2020", "Super Kings won by 5 wickets (with 4 balls remaining)", "Chennai Super Kings", "Mumbai Indians", "MS Dhoni (c)†", "0", "2", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 25 2020", "Super Kings won by 8 wickets (with 8 balls remaining)", "Chennai Super Kings", "Royal Challengers Bangalore", "MS Dhoni (c)†", "19", "21", "3", "0", "90.47"),
    @(" Sharjah", " October 23 2020", "Mumbai won by 10 wickets (with 46 balls remaining)", "Chennai Super Kings", "Mumbai Indians", "MS Dhoni (c)†", "16", "16", "2", "1", "100.00"),
    @(" Dubai (DSC)", " October 02 2020", "Sunrisers won by 7 runs", "Chennai Super Kings", "Sunrisers Hyderabad", "MS Dhoni (c)†", "47", "36", "4", "1", "130.55"),
    @(" Sharjah", " September 22 2020", "Royals won by 16 runs", "Chennai Super Kings", "Rajasthan Royals", "MS Dhoni (c)†", "29", "17", "0", "3", "170.58"),
    @(" Dubai (DSC)", " October 29 2020", "Super Kings won by 6 wickets", "Chennai Super Kings", "Kolkata Knight Riders", "MS Dhoni (c)†", "1", "4", "0", "0", "25.00"),
    @(" Sharjah", " October 17 2020", "Capitals won by 5 wickets (with 1 ball remaining)", "Chennai Super Kings", "Delhi Capitals", "MS Dhoni (c)†", "3", "5", "0", "0", "60.00"),
    @(" Abu Dhabi", " October 19 2020", "Royals won by 7 wickets (with 15 balls remaining)", "Chennai Super Kings", "Rajasthan Royals", "MS Dhoni (c)†", "28", "28", "2", "0", "100.00"),
    @(" Dubai (DSC)", " September 25 2020", "Capitals won by 44 runs", "Chennai Super Kings", "Delhi Capitals", "MS Dhoni (c)†", "15", "12", "2", "0", "125.00"),
    @(" Dubai (DSC)", " October 10 2020", "RCB won by 37 runs", "Chennai Super Kings", "Royal Challengers Bangalore", "MS Dhoni (c)†", "10", "6", "0", "1", "166.66"),
    @(" Abu Dhabi", " October 07 2020", "KKR won by 10 runs", "Chennai Super Kings", "Kolkata Knight Riders", "MS Dhoni (c)†", "11", "12", "1", "0", "91.66")
)

# Columns G-K hold values that look numeric ("21", "161.53", "0.00", ...)
# but must be stored as text (matching the rest of the sheet, which uses
# numberStoredAsText). Setting .Value directly on a numeric-looking string
# makes Excel coerce it to a real number, so we force the cell to a text
# number format first, then restore the default "Normal" style so no extra
# per-cell style is left behind.
$textCols = @(7, 8, 9, 10, 11)

$startRow = 14
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        $colNum = $c + 1
        if ($textCols -contains $colNum) {
            $cell.NumberFormat = "@"
            $cell.Value = $row[$c]
            $cell.Style = "Normal"
        }
        else {
            $cell.Value = $row[$c]
        }
    }
}
